$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Septiembre de 2020 a las 07:25"

# Country name swaps (rank stayed in the same row, but the two countries'
# names traded places because their underlying stats traded places too)
$ws.Range("A63").Value = "Uzbekistan"
$ws.Range("A64").Value = "Nepal"

$ws.Range("A162").Value = "Belice"
$ws.Range("A163").Value = "Birmania"

# Refreshed per-country statistics
$ws.Range("B61").Value = 44199
$ws.Range("C61").Value = 64
$ws.Range("D61").Value = 39389
$ws.Range("E61").Value = 3750
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 1060

$ws.Range("B63").Value = 42903
$ws.Range("C63").Value = 215
$ws.Range("D63").Value = 40176
$ws.Range("E63").Value = 2391
$ws.Range("G63").Value = 5
$ws.Range("H63").Value = 336

$ws.Range("B64").Value = 42877
$ws.Range("D64").Value = 24207
$ws.Range("E64").Value = 18413
$ws.Range("H64").Value = 257

$ws.Range("B73").Value = 26136
$ws.Range("C73").Value = 87
$ws.Range("D73").Value = 22169
$ws.Range("E73").Value = 3230

$ws.Range("B124").Value = 3431
$ws.Range("C124").Value = 4
$ws.Range("E124").Value = 96

$ws.Range("B162").Value = 1118
$ws.Range("D162").Value = 255
$ws.Range("E162").Value = 850
$ws.Range("H162").Value = 13

$ws.Range("B163").Value = 1111
$ws.Range("D163").Value = 359
$ws.Range("E163").Value = 746
$ws.Range("H163").Value = 6

$ws.Range("D187").Value = 150
$ws.Range("E187").Value = 77
